# edit.ps1 - Updates the "as of" date in the confidential disclosure text and
# refreshes the Weight / Percent Change columns (D2:E41) for the QE holdings model.
#
# The worksheet ships protected (no-password "Protect Sheet"), so we briefly
# unprotect it to write the new values and then restore protection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Update the confidential disclosure banner (A44): bump the "as of" date ---
$ws.Range("A44").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# --- Refresh Weight (D) and Percent Change (E) values for each holding row ---
$ws.Range("D2").Value = 0.05862485793293022
$ws.Range("E2").Value = -0.004480212395254313
$ws.Range("D3").Value = 0.05076638047501068
$ws.Range("E3").Value = -0.001603675450763009
$ws.Range("D4").Value = 0.2852690548217338
$ws.Range("E4").Value = -0.03154305200341001
$ws.Range("D5").Value = 0.03506644928793332
$ws.Range("E5").Value = 0.0155119402640036
$ws.Range("D6").Value = 0.0346884530248155
$ws.Range("E6").Value = -0.0159213447510308
$ws.Range("D7").Value = 0.03066494684543358
$ws.Range("E7").Value = -0.06235837940723277
$ws.Range("D8").Value = 0.02973385310907717
$ws.Range("E8").Value = -0.002679628591013894
$ws.Range("D9").Value = 0.02559546571247316
$ws.Range("E9").Value = 0.01330666871779096
$ws.Range("D10").Value = 0.02340866931650076
$ws.Range("E10").Value = 0.002780333837949067
$ws.Range("D11").Value = 0.02226750045107091
$ws.Range("E11").Value = -0.005928853754940788
$ws.Range("D12").Value = 0.0243515754360796
$ws.Range("E12").Value = -0.01052901900359526
$ws.Range("D13").Value = 0.02258643479807658
$ws.Range("E13").Value = 0.04123896346278078
$ws.Range("D14").Value = 0.02125070903251254
$ws.Range("E14").Value = -0.02858310626702998
$ws.Range("D15").Value = 0.0223566723636324
$ws.Range("E15").Value = -0.005361305361305302
$ws.Range("D16").Value = 0.02096408502785059
$ws.Range("E16").Value = 0.009750034525617846
$ws.Range("D17").Value = 0.01967526137214033
$ws.Range("E17").Value = 0.02097753920045209
$ws.Range("D18").Value = 0.01771162637292898
$ws.Range("E18").Value = 0.0005492349941151886
$ws.Range("D19").Value = 0.0151281886775759
$ws.Range("E19").Value = 0.002618041521220027
$ws.Range("D20").Value = 0.01746947251686913
$ws.Range("E20").Value = -0.007338464292106606
$ws.Range("D21").Value = 0.0164475855604501
$ws.Range("E21").Value = -0.005632811124801962
$ws.Range("D22").Value = 0.01630282596213599
$ws.Range("E22").Value = -0.001243118451429548
$ws.Range("D23").Value = 0.01522662520442951
$ws.Range("E23").Value = 0.004745896776745084
$ws.Range("D24").Value = 0.01529170911983153
$ws.Range("E24").Value = -0.008330556481172979
$ws.Range("D25").Value = 0.01367804492550439
$ws.Range("E25").Value = 0.01811870290407258
$ws.Range("D26").Value = 0.01492216681725525
$ws.Range("E26").Value = -0.03967311587623101
$ws.Range("D27").Value = 0.01412853679545793
$ws.Range("E27").Value = 0.003688524590163889
$ws.Range("D28").Value = 0.01326368505129006
$ws.Range("E28").Value = -0.007107184018440571
$ws.Range("D29").Value = 0.01431753492701684
$ws.Range("E29").Value = 0.00310599197618755
$ws.Range("D30").Value = 0.01367005419567746
$ws.Range("E30").Value = -0.0067095331283199
$ws.Range("D31").Value = 0.01313930760441857
$ws.Range("E31").Value = -0.007192088702427157
$ws.Range("D32").Value = 0.01201736281364483
$ws.Range("E32").Value = -0.003372843789149038
$ws.Range("D33").Value = 0.01227839332132485
$ws.Range("E33").Value = 0.01486455896778094
$ws.Range("D34").Value = 0.005893452766564377
$ws.Range("E34").Value = 0.009687561406956258
$ws.Range("D35").Value = 0.005845855810638696
$ws.Range("E35").Value = 0.01463975118366045
$ws.Range("D36").Value = 0.005521941733451025
$ws.Range("E36").Value = 0.0120380856507698
$ws.Range("D37").Value = 0.004851878504774635
$ws.Range("E37").Value = 0.01298453312965431
$ws.Range("D38").Value = 0.005375213404599835
$ws.Range("E38").Value = 0.0222126467736723
$ws.Range("D39").Value = 0.005086041631007552
$ws.Range("E39").Value = 0.005282572066123148
$ws.Range("D40").Value = 0.005162127275881453
$ws.Range("E40").Value = -0.001660123387549062
$ws.Range("E41").Value = -0.0107271494657909

# Restore sheet protection to its original (no password) state.
$ws.Protect()
